# Added effective dates for all transactions to test data
$wb = $excel.ActiveWorkbook

$wsEndorsement   = $wb.Worksheets.Item("EndorsementData")
$wsCancellation  = $wb.Worksheets.Item("CancellationData")
$wsReinstatement = $wb.Worksheets.Item("ReinstatementData")

# EndorsementData!D2 - effective date for the "Add Vehicle" row.
# Leading apostrophe forces text (matches the other hand-typed date-like
# strings already in this workbook, which are stored as quote-prefixed text).
$wsEndorsement.Range("D2").Value = "'12/01/2024"

# CancellationData!B2 - effective date for the cancellation row.
$wsCancellation.Range("B2").Value = "'12/02/2024"
$wsCancellation.Activate()
$wsCancellation.Range("B2").Select()

# ReinstatementData!B2 - effective date for the reinstatement row.
$wsReinstatement.Range("B2").Value = "'19/02/2024"
$wsReinstatement.Activate()
$wsReinstatement.Range("D8").Select()

# EndorsementData!D3 - effective date for the "Add Driver" row.
$wsEndorsement.Range("D3").Value = "'27/01/2024"
$wsEndorsement.Activate()
$wsEndorsement.Range("C7").Select()
